# Mark the TComponentGroup, TLight, TSound and TTimer classes as "Done".
#
# Rows 80-82, 506-510   -> TLight
# Rows 421-425          -> TComponentGroup
# Rows 601-603          -> TSound
# Rows 620-622          -> TTimer
#
# Column layout: A=Class, B=Status, C=Comment, D=Size, E=Function
#
# New shared-strings must come into being in this exact order so they land
# at the same indices as the target workbook (718..721):
#   TComponentGroup, TLight, TSound, TTimer

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- seed the four new class-name strings in the required order ----------
$ws.Range("A421").Value = "TComponentGroup"
$ws.Range("A80").Value = "TLight"
$ws.Range("A601").Value = "TSound"
$ws.Range("A620").Value = "TTimer"

# --- row groups to finish off ---------------------------------------------
$doneRows = @(80, 81, 82, 421, 422, 423, 424, 425, 506, 507, 508, 509, 510, 601, 602, 603, 620, 621, 622)

$classForRow = @{
    80 = "TLight"; 81 = "TLight"; 82 = "TLight";
    421 = "TComponentGroup"; 422 = "TComponentGroup"; 423 = "TComponentGroup"; 424 = "TComponentGroup"; 425 = "TComponentGroup";
    506 = "TLight"; 507 = "TLight"; 508 = "TLight"; 509 = "TLight"; 510 = "TLight";
    601 = "TSound"; 602 = "TSound"; 603 = "TSound";
    620 = "TTimer"; 621 = "TTimer"; 622 = "TTimer"
}

$commentForRow = @{
    421 = "not needed";
    601 = "not needed"
}

foreach ($r in $doneRows) {
    # Copy the "Done" row formatting (fill/border/alignment for A:E) from a
    # known-good reference row (row 4) onto this row, without touching the
    # existing D/E values.
    $ws.Range("A4:E4").Copy()
    $ws.Range("A" + $r + ":E" + $r).PasteSpecial(-4122)

    $ws.Range("A" + $r).Value = $classForRow[$r]
    $ws.Range("B" + $r).Value = "Done"

    if ($commentForRow.ContainsKey($r)) {
        $ws.Range("C" + $r).Value = $commentForRow[$r]
    }
}

# --- selection moved from E2 to D2 in the saved view ----------------------
$ws.Range("D2").Select()
